# Scheduled market-data refresh: update cached price/profit figures
# across several recipe sheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 49
$ws.Range("H49").Value = 499
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null

# Row 76
$ws.Range("H76").Value = 4399.615

# Row 79
$ws.Range("H79").Value = 4399.615

# Row 113
$ws.Range("H113").Value = 3168.875
$ws.Range("I113").Value = 2094.1667
$ws.Range("K113").Value = 2094.1667
$ws.Range("M113").Value = 1159.8333

# Row 137
$ws.Range("H137").Value = 2797.5667
$ws.Range("I137").Value = 2413.625
$ws.Range("J137").Value = 4333.3335
$ws.Range("K137").Value = 7240.875
$ws.Range("L137").Value = 13000.0005
$ws.Range("M137").Value = -4690.875
$ws.Range("N137").Value = -18100.0005

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 80
$ws.Range("H80").Value = 185.66667
$ws.Range("I80").Value = 147.5
$ws.Range("J80").Value = 190.4375
$ws.Range("K80").Value = 147.5
$ws.Range("L80").Value = 190.4375
$ws.Range("M80").Value = 850.5
$ws.Range("N80").Value = -2186.4375

# Row 83
$ws.Range("H83").Value = 185.66667
$ws.Range("I83").Value = 147.5
$ws.Range("J83").Value = 190.4375
$ws.Range("K83").Value = 737.5
$ws.Range("L83").Value = 952.1875
$ws.Range("M83").Value = 4254.5
$ws.Range("N83").Value = -10936.1875

# Row 104
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2866.8125
$ws.Range("I31").Value = 1957.8667
$ws.Range("J31").Value = 3668.8235
$ws.Range("K31").Value = 1957.8667
$ws.Range("L31").Value = 3668.8235
$ws.Range("M31").Value = -1662.8667
$ws.Range("N31").Value = -4258.8235

# Row 34
$ws.Range("H34").Value = 2866.8125
$ws.Range("I34").Value = 1957.8667
$ws.Range("J34").Value = 3668.8235
$ws.Range("K34").Value = 1957.8667
$ws.Range("L34").Value = 3668.8235
$ws.Range("M34").Value = -1755.8667
$ws.Range("N34").Value = -4072.8235

# Row 69
$ws.Range("H69").Value = 23966.666
$ws.Range("I69").Value = 21000
$ws.Range("J69").Value = 24560
$ws.Range("K69").Value = 21000
$ws.Range("L69").Value = 24560
$ws.Range("M69").Value = -20251
$ws.Range("N69").Value = -26058

# Row 72
$ws.Range("H72").Value = 23966.666
$ws.Range("I72").Value = 21000
$ws.Range("J72").Value = 24560
$ws.Range("K72").Value = 63000
$ws.Range("L72").Value = 73680
$ws.Range("M72").Value = -59256
$ws.Range("N72").Value = -81168

# Row 132
$ws.Range("H132").Value = 3108.3
$ws.Range("I132").Value = 2813.05
$ws.Range("J132").Value = 3698.8
$ws.Range("K132").Value = 8439.150000000001
$ws.Range("L132").Value = 11096.4
$ws.Range("M132").Value = -5909.150000000001
$ws.Range("N132").Value = -16156.4

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 55
$ws.Range("H55").Value = 4416.6665
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4416.6665
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = $null
$ws.Range("M55").Value = 13249.9995
$ws.Range("N55").Value = -13603.9995

# Row 74
$ws.Range("H74").Value = 4782
$ws.Range("I74").Value = 350
$ws.Range("J74").Value = 6998
$ws.Range("K74").Value = 1050
$ws.Range("L74").Value = 20994
$ws.Range("M74").Value = 11
$ws.Range("N74").Value = -23116

# Row 77
$ws.Range("H77").Value = 4782
$ws.Range("I77").Value = 350
$ws.Range("J77").Value = 6998
$ws.Range("K77").Value = 3150
$ws.Range("L77").Value = 62982
$ws.Range("M77").Value = 2154
$ws.Range("N77").Value = -73590

# Row 112
$ws.Range("H112").Value = 3005.6843
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 3094.889
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 9284.667000000001
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -11500.667

# Row 122
$ws.Range("H122").Value = 692.3871
$ws.Range("I122").Value = 327.5
$ws.Range("J122").Value = 992.8823
$ws.Range("K122").Value = 2947.5
$ws.Range("L122").Value = 8935.940699999999
$ws.Range("M122").Value = -497.5
$ws.Range("N122").Value = -13835.9407

# Row 125
$ws.Range("H125").Value = 4333.3335
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4333.3335
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = $null
$ws.Range("M125").Value = 13000.0005
$ws.Range("N125").Value = -22840.0005

# Row 130
$ws.Range("H130").Value = 7323.3335
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 7323.3335
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = $null
$ws.Range("M130").Value = 21970.0005
$ws.Range("N130").Value = -32010.0005

# Row 131
$ws.Range("H131").Value = 1304.125
$ws.Range("I131").Value = 2366.25
$ws.Range("J131").Value = 950.0833
$ws.Range("K131").Value = 7098.75
$ws.Range("L131").Value = 2850.2499
$ws.Range("M131").Value = -2058.75
$ws.Range("N131").Value = -12930.2499

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 5747.2
$ws.Range("I70").Value = 5377.8823
$ws.Range("J70").Value = 5971.4287
$ws.Range("K70").Value = 5377.8823
$ws.Range("L70").Value = 5971.4287
$ws.Range("M70").Value = -5107.8823
$ws.Range("N70").Value = -6511.4287

# Row 73
$ws.Range("H73").Value = 5747.2
$ws.Range("I73").Value = 5377.8823
$ws.Range("J73").Value = 5971.4287
$ws.Range("K73").Value = 5377.8823
$ws.Range("L73").Value = 5971.4287
$ws.Range("M73").Value = -4441.8823
$ws.Range("N73").Value = -7843.4287

# Row 80
$ws.Range("H80").Value = 5300.4165
$ws.Range("I80").Value = 2921
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 2921
$ws.Range("L80").Value = 7000
$ws.Range("M80").Value = -1923
$ws.Range("N80").Value = -8996

# Row 83
$ws.Range("H83").Value = 5300.4165
$ws.Range("I83").Value = 2921
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 14605
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = -9613
$ws.Range("N83").Value = -44984

# Row 126
$ws.Range("H126").Value = 2857.1428
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3285.7144
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 9857.143199999999
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -14797.1432

# Row 132
$ws.Range("H132").Value = 2848
$ws.Range("I132").Value = 2186
$ws.Range("K132").Value = 6558
$ws.Range("M132").Value = -4028

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 82
$ws.Range("H82").Value = 1796.6666
$ws.Range("I82").Value = 1767.5
$ws.Range("J82").Value = 1820
$ws.Range("K82").Value = 1767.5
$ws.Range("L82").Value = 1820
$ws.Range("M82").Value = -1406.5
$ws.Range("N82").Value = -2542

# Row 85
$ws.Range("H85").Value = 1796.6666
$ws.Range("I85").Value = 1767.5
$ws.Range("J85").Value = 1820
$ws.Range("K85").Value = 1767.5
$ws.Range("L85").Value = 1820
$ws.Range("M85").Value = -519.5
$ws.Range("N85").Value = -4316

# Row 132
$ws.Range("H132").Value = 5490.0435
$ws.Range("I132").Value = 5376.2856
$ws.Range("J132").Value = 5667
$ws.Range("K132").Value = 16128.8568
$ws.Range("L132").Value = 17001
$ws.Range("M132").Value = -13598.8568
$ws.Range("N132").Value = -22061

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 2087.9092
$ws.Range("I132").Value = 974.8946999999999
$ws.Range("K132").Value = 2924.6841
$ws.Range("M132").Value = -394.6840999999999

# Row 136
$ws.Range("H136").Value = 4399.4
$ws.Range("I136").Value = 4518.4346
$ws.Range("J136").Value = 4274.9546
$ws.Range("K136").Value = 13555.3038
$ws.Range("L136").Value = 12824.8638
$ws.Range("M136").Value = -11005.3038
$ws.Range("N136").Value = -17924.8638

Write-Output "Updated $($wb.Worksheets.Count) worksheets with refreshed market data."
